$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item(4)   # "Iteration #3"
$ws5 = $wb.Worksheets.Item(5)   # "Iteration #4"

# ---------------------------------------------------------------------------
# Sheet "Iteration #3" (sheet4.xml): fill in row 18, mark rows 19-20 as
# "date" formatted (copy format from an already-date-formatted cell) even
# though they stay empty.
# ---------------------------------------------------------------------------
$ws4.Range("A18").Value = 42829
$ws4.Range("A15").Copy()
$ws4.Range("A18").PasteSpecial(-4122)

$ws4.Range("B18").Value = "travail"
$ws4.Range("C18").Value = 2

$ws4.Range("A15").Copy()
$ws4.Range("A19").PasteSpecial(-4122)

$ws4.Range("A15").Copy()
$ws4.Range("A20").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "Iteration #4" (sheet5.xml): fill in the work log, rows 14-21.
# The shared-strings table is built in the order new unique strings are
# written, so rows 20/21 are populated before rows 18/19 to land on the same
# shared-string indices the original author ended up with.
# ---------------------------------------------------------------------------
$ws5.Range("A14").Value = 42835
$ws5.Range("B14").Value = "travail"
$ws5.Range("C14").Value = 3

$ws5.Range("A15").Value = 42836
$ws5.Range("B15").Value = "travail"
$ws5.Range("C15").Value = 2

$ws5.Range("A16").Value = 42842
$ws5.Range("B16").Value = "travail"
$ws5.Range("C16").Value = 3

$ws5.Range("A17").Value = 42843
$ws5.Range("B17").Value = "travail"
$ws5.Range("C17").Value = 2

$ws5.Range("A20").Value = 42856
$ws5.Range("B20").Value = "implémenter les commandes sonores dans l'appli"
$ws5.Range("C20").Value = 3

$ws5.Range("A21").Value = 42857
$ws5.Range("B21").Value = "finaliser d'implémenter les commandes sonores dans l'appli"
$ws5.Range("C21").Value = 2

$ws5.Range("A18").Value = 42849
$ws5.Range("B18").Value = "apprendre à implémenter des commandes snonres dans une appli"
$ws5.Range("C18").Value = 3

$ws5.Range("A19").Value = 42850
$ws5.Range("B19").Value = "création des sonds à implémenter dans l'ap^pli"
$ws5.Range("C19").Value = 2

# Give the date cells A15:A21 the same "date" number format as Iteration #3's
# A15 (style 23), so they match the newly-dated rows without touching A14
# (which keeps its original style 12).
$ws4.Range("A15").Copy()
$ws5.Range("A15:A21").PasteSpecial(-4122)

$wb.Application.Calculate()

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: the author left off on "Iteration #3"
# with C18 selected, then moved on to "Iteration #4" with B18 selected as the
# final active sheet.
# ---------------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("C18").Select()

$ws5.Activate()
$ws5.Range("B18").Select()
